function Set-CellText($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "42.850.09"
Set-CellText $ws "E2" "  +0.11%  "
Set-CellText $ws "D3" "2.366.79"
Set-CellText $ws "E3" "  +2.12%  "
Set-CellText $ws "E4" "  +0.08%  "
Set-CellText $ws "D5" "301.31"
Set-CellText $ws "E5" "  -0.41%  "
Set-CellText $ws "E6" "  -0.02%  "
Set-CellText $ws "E7" "  -0.50%  "
Set-CellText $ws "E8" "  -0.03%  "
Set-CellText $ws "D9" "0.492"
Set-CellText $ws "E9" "  -0.44%  "
Set-CellText $ws "D10" "34.02"
Set-CellText $ws "E10" "  -1.31%  "
Set-CellText $ws "E11" "  +0.23%  "
Set-CellText $ws "E12" "  +2.48%  "
Set-CellText $ws "D13" "18.29"
Set-CellText $ws "E13" "  -3.27%  "
Set-CellText $ws "D14" "6.75"
Set-CellText $ws "E14" "  -0.02%  "
Set-CellText $ws "D15" "2.739.95"
Set-CellText $ws "E15" "  +2.28%  "
Set-CellText $ws "D16" "2.354.74"
Set-CellText $ws "E16" "  +1.48%  "
Set-CellText $ws "D17" "0.799"
Set-CellText $ws "E17" "  +0.76%  "
Set-CellText $ws "D18" "42.828.65"
Set-CellText $ws "E18" "  +0.19%  "
Set-CellText $ws "D19" "12.14"
Set-CellText $ws "E19" "  -0.15%  "
Set-CellText $ws "D20" "6.29"
Set-CellText $ws "E20" "  +2.05%  "
Set-CellText $ws "E21" "  -0.72%  "
Set-CellText $ws "D22" "67.95"
Set-CellText $ws "E22" "  +0.12%  "
Set-CellText $ws "D23" "234.82"
Set-CellText $ws "E23" "  -0.48%  "
Set-CellText $ws "E24" "  -1.68%  "
Set-CellText $ws "E25" "  -0.06%  "
Set-CellText $ws "E26" "  +0.41%  "
Set-CellText $ws "D27" "24.77"
Set-CellText $ws "E27" "  +1.80%  "
Set-CellText $ws "D28" "2.37"
Set-CellText $ws "E28" "  -0.27%  "
Set-CellText $ws "D29" "9.19"
Set-CellText $ws "E29" "  +0.53%  "
Set-CellText $ws "D30" "31.49"
Set-CellText $ws "E30" "  -2.42%  "
Set-CellText $ws "E31" "  +0.09%  "
Set-CellText $ws "E32" "  +0.61%  "
Set-CellText $ws "D33" "0.0734"
Set-CellText $ws "E33" "  +4.97%  "
Set-CellText $ws "D34" "17.20"
Set-CellText $ws "E34" "  -4.08%  "
Set-CellText $ws "E35" "  +5.30%  "
Set-CellText $ws "E36" "  +4.21%  "
Set-CellText $ws "E37" "  -2.30%  "
Set-CellText $ws "D38" "2.30"
Set-CellText $ws "E38" "  -1.34%  "
Set-CellText $ws "E39" "  +1.63%  "
Set-CellText $ws "D40" "22.09"
Set-CellText $ws "E40" "  +6.39%  "
Set-CellText $ws "E41" "  -0.63%  "
Set-CellText $ws "D42" "117.56"
Set-CellText $ws "E42" "  -29.31%  "
Set-CellText $ws "D43" "1.936.19"
Set-CellText $ws "E43" "  +0.41%  "
Set-CellText $ws "E44" "  +0.32%  "
Set-CellText $ws "E46" "  -0.93%  "
Set-CellText $ws "E47" "  -9.50%  "
Set-CellText $ws "D48" "2.597.81"
Set-CellText $ws "E48" "  +2.01%  "
Set-CellText $ws "E49" "  +1.86%  "
Set-CellText $ws "B50" "BitcoinSV"
Set-CellText $ws "C50" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-CellText $ws "D50" "71.96"
Set-CellText $ws "E50" "  -0.18%  "
Set-CellText $ws "B51" "MultiversX"
Set-CellText $ws "C51" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-CellText $ws "D51" "51.94"
Set-CellText $ws "E51" "  -2.63%  "
